$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A4: "Sanitary facilities - Access to toilet" -> "Sanitary facilities – No Access to toilet"
$ws.Range("A4").Value = "Sanitary facilities – No Access to toilet"

# --- A5: "Access to toilet" -> "crash test"
$ws.Range("A5").Value = "crash test"

# --- B4 rich text: "SELECT COUNT(*) FROM reg_entities;" -> "SELECT COUNT(*) FROM reg_entity_farmaq_details;"
# (last two runs " " + "reg_entities;" become " reg_entity_farmaq_details" + ";")
$cell = $ws.Range("B4")
$cell.Value = "SELECT COUNT(*) FROM reg_entity_farmaq_details;"

# Run 1: "SELECT" (bold orange)
$cell.Characters(1,6).Font.Bold = $true
$cell.Characters(1,6).Font.Color = 0x4080FF
$cell.Characters(1,6).Font.Name = "Monospace"
$cell.Characters(1,6).Font.Size = 10

# Run 2: " " (gray)
$cell.Characters(7,1).Font.Bold = $false
$cell.Characters(7,1).Font.Color = 0x969483
$cell.Characters(7,1).Font.Name = "Monospace"
$cell.Characters(7,1).Font.Size = 10

# Run 3: "COUNT" (bold yellow)
$cell.Characters(8,5).Font.Bold = $true
$cell.Characters(8,5).Font.Color = 0x80FFFF
$cell.Characters(8,5).Font.Name = "Monospace"
$cell.Characters(8,5).Font.Size = 10

# Run 4: "(*) " (light gray)
$cell.Characters(13,4).Font.Bold = $false
$cell.Characters(13,4).Font.Color = 0xDDDDDD
$cell.Characters(13,4).Font.Name = "Monospace"
$cell.Characters(13,4).Font.Size = 10

# Run 5: "FROM" (bold orange)
$cell.Characters(17,4).Font.Bold = $true
$cell.Characters(17,4).Font.Color = 0x4080FF
$cell.Characters(17,4).Font.Name = "Monospace"
$cell.Characters(17,4).Font.Size = 10

# Run 6: " reg_entity_farmaq_details" (gray)
$cell.Characters(21,26).Font.Bold = $false
$cell.Characters(21,26).Font.Color = 0x969483
$cell.Characters(21,26).Font.Name = "Monospace"
$cell.Characters(21,26).Font.Size = 10

# Run 7: ";" (light gray)
$cell.Characters(47,1).Font.Bold = $false
$cell.Characters(47,1).Font.Color = 0xDDDDDD
$cell.Characters(47,1).Font.Name = "Monospace"
$cell.Characters(47,1).Font.Size = 10

# --- sheet view: selection moves from B10 to B4
[void]$ws.Range("B4").Select()

# --- column widths: A 24.72 -> 52.59, B 29 -> 64.17
# (closest achievable values given Excel's column-width pixel quantization)
$ws.Columns.Item(1).ColumnWidth = 51.83
$ws.Columns.Item(2).ColumnWidth = 63.33
